# 18 Mayıs 2020 verileri eklendi
# Adds the 2020-05-18 COVID-19 Turkey daily row to the "data" sheet/table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 68): date, test, case, death, recovered
$ws.Range("A68").Value = 43969
$ws.Range("B68").Value = 25141
$ws.Range("C68").Value = 1158
$ws.Range("D68").Value = 31
$ws.Range("E68").Value = 1615

# Grow the Excel table ("Table3") so the new row is included, which also
# extends the AutoFilter range to match.
$tbl = $ws.ListObjects("Table3")
$tbl.Resize($ws.Range("A1:E68"))

# Mirror the author's last active cell after entering the new row.
[void]$ws.Range("E67").Select()
